# ---------------------------------------------------------------------------
# Update "Data Attribute Description File.xlsx":
#   - rename Sheet1 -> "Attributes Initial Load"
#   - rewrite column C (description) text for several rows to reflect the
#     newer Olist data-dictionary wording
#   - add a new "Missing Data" column D with per-attribute missing-data stats
#   - clear out the two trailing rows that used to hold the now-removed
#     purchase_wk_day / purchase_month attributes
#   - tidy up column widths / header formatting / active-sheet & selection
# ---------------------------------------------------------------------------

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# 1. Rename the first sheet.
$ws1.Name = "Attributes Initial Load"

$cData = @{}
$dData = @{}
$cData[2] = 'order unique identifier'
$dData[2] = 0
$cData[3] = 'key to the orders dataset - each order has a unique customer_id'
$dData[3] = 0
$cData[4] = 'order status, 7-levels (shipped, canceled, invoiced, processing, approved, unavailable, delivered)'
$dData[4] = 0
$cData[5] = 'purchase initiation timestamp'
$dData[5] = 0
$cData[6] = 'payment approval timestamp '
$dData[6] = 160
$cData[7] = 'order posting timestamp when it was handed to the logistic partner '
$dData[7] = ' 2% = 1,783'
$cData[8] = 'actual order delivery date to the customer '
$dData[8] = '3% = 2,965'
$cData[9] = 'estimated delivery date provided to the customer at the time of purchase initiation'
$dData[9] = 0
$cData[10] = 'unique identifier of a customer'
$dData[10] = 0
$cData[11] = 'first five digits of customer zip code'
$dData[11] = 0
$cData[12] = 'customer city name'
$dData[12] = 0
$cData[13] = 'customer state name'
$dData[13] = 0
$cData[14] = 'sequential number identifying number of items included in the same order'
$dData[14] = 0
$cData[15] = 'product unique identifier'
$dData[15] = 0
$cData[16] = 'seller unique identifier'
$dData[16] = 0
$cData[17] = 'seller shipping limit date for handing the order off to the logistic partner'
$dData[17] = 0
$cData[18] = 'item price'
$dData[18] = 0
$cData[19] = 'item freight value (if an order has more than one item, the freight value is split between the items)'
$dData[19] = 0
$cData[20] = 'number of payment methods used by the customer'
$dData[20] = 0
$cData[21] = 'method of payment by customer [74% credit_card, 19% boleto, 7% other]'
$dData[21] = 0
$cData[22] = 'number of payment installments by customer'
$dData[22] = 0
$cData[23] = 'transaction value'
$dData[23] = 0
$cData[24] = 'first five digits of seller zip code'
$dData[24] = 0
$cData[25] = 'seller city name'
$dData[25] = 0
$cData[26] = 'seller state name'
$dData[26] = 0
$cData[27] = 'root category of product in Portuguese'
$dData[27] = '2% = 610 '
$cData[28] = 'number of characters extracted from the product name'
$dData[28] = '2% = 610 '
$cData[29] = 'number of characters extracted from the product description'
$dData[29] = '2% = 610 '
$cData[30] = 'number of product photos published'
$dData[30] = '2% = 610 '
$cData[31] = 'product weight measured in grams'
$dData[31] = 2
$cData[32] = 'product length measured in centimeters'
$dData[32] = 2
$cData[33] = 'product height measured in cemitmeters'
$dData[33] = 2
$cData[34] = 'product width measured in centimeters'
$dData[34] = 2
$cData[35] = 'product category name in English'
$dData[35] = '2% = 610 '
$cData[36] = 'review unique identifier'
$dData[36] = 0
$cData[37] = '1 to 5 rating given by the customer on a satisfaction survey'
$dData[37] = 0
$cData[38] = 'comment titles from the review left by the customer'
$dData[38] = '88% = 88.3K'
$cData[39] = 'comment message from the review left by the customer [note: 58% missing]'
$dData[39] = '58%=58.2k'
$cData[40] = 'date satisfaction survey sent to customer'
$dData[40] = 0
$cData[41] = 'satisfaction survey answer timestamp'
$dData[41] = 0

# 2. Rewrite column C description text for rows 2-41, and populate the new
#    column D ("Missing Data") values for the same rows.
foreach ($r in $cData.Keys) {
    $ws1.Cells.Item($r, 3).Value = $cData[$r]
    $ws1.Cells.Item($r, 4).Value = $dData[$r]
}

# 3. Rows 42 & 43 used to describe purchase_wk_day / purchase_month, which no
#    longer exist as attributes - blank out their contents (formatting stays).
$ws1.Cells.Item(42, 1).Value = ""
$ws1.Cells.Item(42, 2).Value = ""
$ws1.Cells.Item(42, 3).Value = ""
$ws1.Cells.Item(43, 1).Value = ""
$ws1.Cells.Item(43, 2).Value = ""
$ws1.Cells.Item(43, 3).Value = ""

# 4. Header row (row 1) becomes bold; new D1 header ("Missing Data") is bold
#    and center-aligned, matching the font size used by the rest of the
#    header (10pt) rather than the workbook default (11pt).
$ws1.Range("A1:C1").Font.Bold = $true
$ws1.Cells.Item(1, 4).Value = "Missing Data"
$ws1.Range("D1").Font.Size = 10
$ws1.Range("D1").Font.Bold = $true
$ws1.Range("D1").HorizontalAlignment = -4108

# 5. Data cells in the new column D (rows 2-41) are center-aligned.
$ws1.Range("D2:D41").HorizontalAlignment = -4108

# 6. Column width tweaks: description column narrows, new column gets a
#    sensible width.
$ws1.Columns.Item(3).ColumnWidth = 78.8
$ws1.Columns.Item(4).ColumnWidth = 10.6

# 7. View state: "Attributes Initial Load" becomes the active/selected sheet
#    (with cell C13 selected, scrolled back to the top), and the
#    "Missing Data Pass 1" sheet loses its stale tabSelected flag.
$ws1.Activate() | Out-Null
$ws1.Range("C13").Select() | Out-Null

Write-Host "done"
